# Record the PAM fluorometry species for every measurement row.
#
# The source data only ever covers a single cross: the aggregating
# anemone A. elegantissima x B. muscatinei. Add a new "species" column
# (N) with that label, header in N1, and a value for every data row
# (rows 3-62; row 2 is a sub-header row and stays blank in column N,
# matching the rest of the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$headerText = "species"
$speciesText = "A. elegantissima x B. muscatinei"

# Header for the new column.
$ws.Range("N1").Value = $headerText

# Fill the species value down every data row (3 through 62).
$ws.Range("N3:N62").Value = $speciesText

# Mirror the author's final selection/view state after entering the data.
$ws.Range("N3:N62").Select()
